# Duplicate the "CCC" sheet, placing the copy immediately after it and
# naming it "CCC_Dup" (sheetId=4, rId4) - matches the commit:
# "add deleted, new, renamed sheets to output".

$wb = $excel.ActiveWorkbook

$ccc = $wb.Worksheets.Item("CCC")

# Copy CCC to a new sheet placed right after CCC.
$ccc.Copy($null, $ccc)

# The freshly copied sheet is the last one / currently active; rename it.
$dup = $wb.Worksheets.Item($wb.Worksheets.Count)
$dup.Name = "CCC_Dup"

# Keep "CCC" as the selected/active tab (matches activeTab pointing at CCC).
$ccc.Activate()
